$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.380099999999996
$ws.Range("B3").Value = 5.787799999999992
$ws.Range("B5").Value = 4.890400000000001
$ws.Range("D7").Value = -7.286600000000002
$ws.Range("A9").Value = -20.28159999999997
$ws.Range("D9").Value = -8.617400000000005
$ws.Range("B11").Value = 5.444599999999997
$ws.Range("B12").Value = 5.410799999999998
$ws.Range("A13").Value = -21.91350000000003
$ws.Range("A16").Value = -20.10429999999999
$ws.Range("A18").Value = -22.5012
$ws.Range("A20").Value = -21.98310000000001
$ws.Range("B21").Value = 5.492199999999992
$ws.Range("D21").Value = -7.692599999999995
